$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update "Marking" row (row 11): Right marks per question 3 -> 5
$ws.Range("B11").Value = 5

# Update "Total" row (row 12): Right total 60 -> 100
$ws.Range("B12").Value = 100

# Update Correct/Total marks label 52/84 -> 100/140
$ws.Range("E12").Value = "100/140"
